$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A6").Value = "final"
$ws.Rows.Item(6).RowHeight = 21.75
